$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the last "Brasil" row (row 27), then after the last
# "Nordeste" row (now row 54), then append a new row after the last
# "Sergipe" row (now row 81) - each carrying the new 01/07/2025 quarter.

$varLabel = "Taxa de pessoas de 14 anos ou mais de idade, na força de trabalho, na semana de referência"

# 1) Insert new Brasil row after row 27 (shifts Nordeste/Sergipe blocks down by 1)
$ws.Rows.Item(28).Insert()
$ws.Cells.Item(28, 1).Value = "Brasil"
$ws.Cells.Item(28, 2).Value = $varLabel
$ws.Cells.Item(28, 3).NumberFormat = "@"
$ws.Cells.Item(28, 3).Value = "01/07/2025"
$ws.Cells.Item(28, 3).ClearFormats()
$ws.Cells.Item(28, 4).Value = 50.95

# 2) Insert new Nordeste row after the (now shifted) last Nordeste row, row 54
$ws.Rows.Item(55).Insert()
$ws.Cells.Item(55, 1).Value = "Nordeste"
$ws.Cells.Item(55, 2).Value = $varLabel
$ws.Cells.Item(55, 3).NumberFormat = "@"
$ws.Cells.Item(55, 3).Value = "01/07/2025"
$ws.Cells.Item(55, 3).ClearFormats()
$ws.Cells.Item(55, 4).Value = 44.31

# 3) Insert new Sergipe row after the (now shifted) last Sergipe row, row 81
$ws.Rows.Item(82).Insert()
$ws.Cells.Item(82, 1).Value = "Sergipe"
$ws.Cells.Item(82, 2).Value = $varLabel
$ws.Cells.Item(82, 3).NumberFormat = "@"
$ws.Cells.Item(82, 3).Value = "01/07/2025"
$ws.Cells.Item(82, 3).ClearFormats()
$ws.Cells.Item(82, 4).Value = 43.22
